# Update "Datos actualizados" timestamp in A1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 1 de Septiembre de 2020 a las 11:24"

# Country ranking swaps (column A) caused by updated case counts
$ws.Range("A47").Value = "Polonia"
$ws.Range("A48").Value = "Japon"

$ws.Range("A72").Value = "El Salvador"
$ws.Range("A73").Value = "Australia"

$ws.Range("A119").Value = "Eslovaquia"
$ws.Range("A120").Value = "Congo"

# Updated numeric statistics (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes)

# Row 25
$ws.Range("B25").Value = 224264
$ws.Range("C25").Value = 3483
$ws.Range("D25").Value = 158012
$ws.Range("E25").Value = 62655
$ws.Range("G25").Value = 39
$ws.Range("H25").Value = 3597

# Row 26
$ws.Range("B26").Value = 177571
$ws.Range("C26").Value = 2775
$ws.Range("D26").Value = 128057
$ws.Range("E26").Value = 42009
$ws.Range("G26").Value = 88
$ws.Range("H26").Value = 7505

# Row 30
$ws.Range("B30").Value = 117241
$ws.Range("C30").Value = 645
$ws.Range("D30").Value = 95596
$ws.Range("E30").Value = 20699
$ws.Range("G30").Value = 7
$ws.Range("H30").Value = 946

# Row 47
$ws.Range("B47").Value = 67922
$ws.Range("C47").Value = 550
$ws.Range("D47").Value = 47030
$ws.Range("E47").Value = 18834
$ws.Range("G47").Value = 19
$ws.Range("H47").Value = 2058

# Row 48
$ws.Range("B48").Value = 67865
$ws.Range("D48").Value = 56802
$ws.Range("E48").Value = 9784
$ws.Range("H48").Value = 1279

# Row 65
$ws.Range("B65").Value = 38196
$ws.Range("C65").Value = 31
$ws.Range("D65").Value = 29231
$ws.Range("E65").Value = 7559
$ws.Range("G65").Value = 4
$ws.Range("H65").Value = 1406

# Row 71
$ws.Range("B71").Value = 27642
$ws.Range("C71").Value = 204
$ws.Range("D71").Value = 23565
$ws.Range("E71").Value = 3343
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 734

# Row 72
$ws.Range("B72").Value = 25820
$ws.Range("C72").Value = 91
$ws.Range("D72").Value = 14570
$ws.Range("E72").Value = 10526
$ws.Range("G72").Value = 7
$ws.Range("H72").Value = 724

# Row 73
$ws.Range("B73").Value = 25819
$ws.Range("C73").Value = 73
$ws.Range("D73").Value = 21503
$ws.Range("E73").Value = 3659
$ws.Range("G73").Value = 5
$ws.Range("H73").Value = 657

# Row 101
$ws.Range("B101").Value = 8142
$ws.Range("C101").Value = 56
$ws.Range("E101").Value = 606

# Row 112
$ws.Range("B112").Value = 4823
$ws.Range("C112").Value = 12
$ws.Range("D112").Value = 4380
$ws.Range("E112").Value = 353
$ws.Range("G112").Value = 1
$ws.Range("H112").Value = 90

# Row 119
$ws.Range("B119").Value = 3989
$ws.Range("C119").Value = 72
$ws.Range("D119").Value = 2478
$ws.Range("E119").Value = 1478
$ws.Range("H119").Value = 33

# Row 120
$ws.Range("B120").Value = 3979
$ws.Range("D120").Value = 1742
$ws.Range("E120").Value = 2159
$ws.Range("H120").Value = 78

# Row 131
$ws.Range("B131").Value = 2924
$ws.Range("C131").Value = 41
$ws.Range("D131").Value = 2323
$ws.Range("E131").Value = 468

# Row 137
$ws.Range("B137").Value = 2395
$ws.Range("C137").Value = 20
$ws.Range("D137").Value = 2112
$ws.Range("E137").Value = 219

# Row 152
$ws.Range("D152").Value = 1253
$ws.Range("E152").Value = 238
